$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. LinqWhere sheet: insert the new "For Loop" benchmark row, append a new
#    trailing row, and refresh the benchmark numbers that shifted because of
#    the newly measured scenario.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LinqWhere")

# Make room for the new "GetUsingMergeSetsAndSingleToUpper" benchmark just
# after the existing row 3 (shifts the old rows 4-8 down to 5-9, preserving
# their per-cell number formatting).
$ws.Rows(4).Insert()

$ws.Range("A4").Value = "GetUsingMergeSetsAndSingleToUpper"
$ws.Range("B4").Value = 624.6
$ws.Range("B4").NumberFormat = "#,##0.00"
$ws.Range("C4").Value = 1.08
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 0.0906
$ws.Range("F4").Value = 384

# Append the new "GetWithoutLINQNoToUpperForLoop" benchmark as the new last
# row of the table.
$ws.Range("A10").Value = "GetWithoutLINQNoToUpperForLoop"
$ws.Range("B10").Value = 165.1
$ws.Range("C10").Value = 0.28
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.0093
$ws.Range("F10").Value = 40

# Refresh every benchmark row with the re-measured results.
$ws.Range("B2").Value = 579.9
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0.103
$ws.Range("F2").Value = 432

$ws.Range("B3").Value = 1000.1
$ws.Range("C3").Value = 1.72
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 0.1144
$ws.Range("F3").Value = 488

$ws.Range("B5").Value = 1928.1
$ws.Range("C5").Value = 3.32
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 0.2213
$ws.Range("F5").Value = 944

$ws.Range("B6").Value = 691.9
$ws.Range("C6").Value = 1.19
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 0.123
$ws.Range("F6").Value = 520

$ws.Range("B7").Value = 431.2
$ws.Range("C7").Value = 0.74
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 0.0725
$ws.Range("F7").Value = 304

$ws.Range("B8").Value = 288.2
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 0.0453
$ws.Range("F8").Value = 192

$ws.Range("B9").Value = 206
$ws.Range("C9").Value = 0.36
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 0.0188
$ws.Range("F9").Value = 80

# Grow Table1 so it covers the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F10"))

# Add the bar chart that visualises the LinqWhere benchmarks.
$chartObj = $ws.ChartObjects().Add(260, 0, 620, 400)
$chartObj.Chart.ChartType = 51
$chartObj.Chart.SetSourceData($ws.Range("A1:B10"))
$chartObj.Chart.HasTitle = $true
$chartObj.Chart.ChartTitle.Text = "To LINQ or Not to LINQ"

# ---------------------------------------------------------------------------
# 2. SlugProducer sheet: selection moved from K15 to A3.
# ---------------------------------------------------------------------------
$wsSlug = $wb.Worksheets.Item("SlugProducer")
$wsSlug.Activate()
$wsSlug.Range("A3").Select()

# ---------------------------------------------------------------------------
# 3. Re-select C15 on LinqWhere and make it the active (visible) sheet/tab,
#    matching the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C15").Select()
